$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A. This shifts A->B, B->C, D->E, E->F
# (the old, always-empty column C simply slides into the new D slot).
$ws.Columns.Item(1).Insert()

# Fill the two "index" columns (A and D) with row numbers 1-12
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
    $ws.Cells.Item($r, 4).Value = $r - 1
}

# New "Correct Key" block (added first so its strings land earlier in the
# shared-string table, matching the authored workbook)
$ws.Range("B16").Value = "Correct Key:"
$ws.Range("C16").Value = "110100100"

# New header cell for the "correct output" column
$ws.Range("C1").Value = "correct output"

# Mark G2 as a text-formatted (but otherwise empty) cell
$ws.Range("G2").NumberFormat = "@"

# Give the new trailing column a sensible custom width (close to the
# authored workbook's auto-fit result)
$ws.Columns.Item(7).ColumnWidth = 26.5

# Update selection to match the authored workbook
$ws.Range("G2").Select()
